# Update "想去人数" (number of people interested) figures across sheets
# to reflect the latest scrape, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5419
$ws1.Range("F6").Value = 816
$ws1.Range("F7").Value = 15
$ws1.Range("F8").Value = 338
$ws1.Range("F9").Value = 15

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 15

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5419
$ws4.Range("F6").Value = 816
$ws4.Range("F7").Value = 15
$ws4.Range("F9").Value = 338
$ws4.Range("F10").Value = 15
$ws4.Range("F11").Value = 15
